$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# Both locale sheets (zh-cn, de-de) get two new columns of data filled in
# (Latest Target File / Latest Handback File), the "Latest Handback
# DateTime" column gets a real timestamp instead of the epoch placeholder,
# and the Status column moves from "Ready for handoff" to a handed-back
# message.
# ---------------------------------------------------------------------------

function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$HandoffXlf1,
        [string]$HandoffXlf2,
        [string]$HandbackTime1,
        [string]$HandbackTime2
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $status = "Handed back: in sync with en-US"

    # Row 2 - 645e67dd-...md
    $ws.Range("B2").Value = $status
    $ws.Range("E2").Value = $ws.Range("A2").Value2
    $ws.Range("F2").Value = $HandoffXlf1
    $ws.Range("G2").Value = $HandbackTime1

    # Row 3 - c7a2c94e-...md
    $ws.Range("B3").Value = $status
    $ws.Range("E3").Value = $ws.Range("A3").Value2
    $ws.Range("F3").Value = $HandoffXlf2
    $ws.Range("G3").Value = $HandbackTime2

    # Recreate every hyperlink on the sheet so the new Latest Target
    # File / Latest Handback File columns get a link too, in the same
    # left-to-right, top-to-bottom order Excel would naturally emit them.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $ws.Range("A2").Value2, [Type]::Missing, [Type]::Missing, $ws.Range("A2").Value2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $ws.Range("C2").Value2, [Type]::Missing, [Type]::Missing, $ws.Range("C2").Value2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E2"), $ws.Range("E2").Value2, [Type]::Missing, [Type]::Missing, $ws.Range("E2").Value2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $ws.Range("F2").Value2, [Type]::Missing, [Type]::Missing, $ws.Range("F2").Value2) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), $ws.Range("A3").Value2, [Type]::Missing, [Type]::Missing, $ws.Range("A3").Value2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $ws.Range("C3").Value2, [Type]::Missing, [Type]::Missing, $ws.Range("C3").Value2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E3"), $ws.Range("E3").Value2, [Type]::Missing, [Type]::Missing, $ws.Range("E3").Value2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $ws.Range("F3").Value2, [Type]::Missing, [Type]::Missing, $ws.Range("F3").Value2) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), $ws.Range("A4").Value2, [Type]::Missing, [Type]::Missing, $ws.Range("A4").Value2) | Out-Null
}

Update-LocaleSheet "zh-cn" `
    "645e67dd-0c23-4a3b-91bb-bfc452660f53.c0c72e6ce361e443ed4869cfcecab95ca4268589.zh-cn.xlf" `
    "c7a2c94e-e6af-4aa1-a2e0-c1e3208c0a73.44d208c53bac163d416128242ff2c461a4baf71a.zh-cn.xlf" `
    "2016-02-18 10:30:47" `
    "2016-02-18 10:30:47"

Update-LocaleSheet "de-de" `
    "645e67dd-0c23-4a3b-91bb-bfc452660f53.c0c72e6ce361e443ed4869cfcecab95ca4268589.de-de.xlf" `
    "c7a2c94e-e6af-4aa1-a2e0-c1e3208c0a73.44d208c53bac163d416128242ff2c461a4baf71a.de-de.xlf" `
    "2016-02-18 10:31:10" `
    "2016-02-18 10:31:10"
